# Complete list chooser added
# The underlying "Количество элементов" (element count) samples were
# regenerated with 15 elements instead of 100, which cascades into new
# computed statistics (geometric mean, arithmetic mean, std-dev estimate,
# covariance coefficient, coefficient of variation, confidence-interval
# bounds, variance estimate, max/min) for all three samples. Column F's
# best-fit width also shrank slightly to match its new (shorter) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sample 1 (columns A/B) ---
$ws.Range("B2").Value  = 0.0019104679425557306
$ws.Range("B3").Value  = 0.5892285904932517
$ws.Range("B4").Value  = 1.6918176305480301
$ws.Range("B5").Value  = -0.0016375077264237116
$ws.Range("B6").Value  = 15.0
$ws.Range("B7").Value  = 297.96307331355297
$ws.Range("B8").Value  = -0.24911750904730395
$ws.Range("B9").Value  = 0.2529384449324154
$ws.Range("B10").Value = 0.3471903318546641
$ws.Range("B11").Value = 0.8207750557921827
$ws.Range("B12").Value = -0.8710425747558475

# --- Sample 2 (columns C/D) ---
$ws.Range("D1").Value  = 0.5565811503729629
$ws.Range("D2").Value  = 0.6151980836500521
$ws.Range("D3").Value  = 0.23824862442915698
$ws.Range("D4").Value  = 0.8246188284448921
$ws.Range("D5").Value  = 0.05067571442929425
$ws.Range("D6").Value  = 15.0
$ws.Range("D7").Value  = 0.3741397193978958
$ws.Range("D8").Value  = 0.5136974586362476
$ws.Range("D9").Value  = 0.7166987086638567
$ws.Range("D10").Value = 0.05676240704238549
$ws.Range("D11").Value = 0.9662961962314939
$ws.Range("D12").Value = 0.14167736778660175

# --- Sample 3 (columns E/F) ---
$ws.Range("F2").Value  = 1.3706230941087303
$ws.Range("F3").Value  = 3.4716991967473048
$ws.Range("F4").Value  = 11.13785641232145
$ws.Range("F5").Value  = 0.0053383090152330295
$ws.Range("F6").Value  = 15.0
$ws.Range("F7").Value  = 2.4470476834344765
$ws.Range("F8").Value  = -0.10841855481610696
$ws.Range("F9").Value  = 2.8496647430335678
$ws.Range("F10").Value = 12.05269531269588
$ws.Range("F11").Value = 7.439512783585943
$ws.Range("F12").Value = -3.698343628735506

# --- Column F best-fit width shrank (13.8515625 -> 12.73828125) ---
$ws.Columns.Item(6).ColumnWidth = 11.833333333333334
